$d = $word.ActiveDocument
$d.Content.Find.Execute("22+57=", $true, $false, $false, $false, $false, $true, 1, $false, "80-51=", 2)
$d.Content.Find.Execute("10+6=", $true, $false, $false, $false, $false, $true, 1, $false, "77-14=", 2)
$d.Content.Find.Execute("90-23=", $true, $false, $false, $false, $false, $true, 1, $false, "35+54=", 2)
$d.Content.Find.Execute("86-65=", $true, $false, $false, $false, $false, $true, 1, $false, "44-26=", 2)
$d.Content.Find.Execute("12+66=", $true, $false, $false, $false, $false, $true, 1, $false, "30-16=", 2)
$d.Content.Find.Execute("91+5=", $true, $false, $false, $false, $false, $true, 1, $false, "59+27=", 2)
$d.Content.Find.Execute("48+37=", $true, $false, $false, $false, $false, $true, 1, $false, "11+77=", 2)
$d.Content.Find.Execute("34+54=", $true, $false, $false, $false, $false, $true, 1, $false, "23-11=", 2)
$d.Content.Find.Execute("36+47=", $true, $false, $false, $false, $false, $true, 1, $false, "99-74=", 2)
$d.Content.Find.Execute("53+46=", $true, $false, $false, $false, $false, $true, 1, $false, "45-29=", 2)
$d.Content.Find.Execute("43+53=", $true, $false, $false, $false, $false, $true, 1, $false, "90+9=", 2)
$d.Content.Find.Execute("89+0=", $true, $false, $false, $false, $false, $true, 1, $false, "31-31=", 2)
$d.Content.Find.Execute("43+55=", $true, $false, $false, $false, $false, $true, 1, $false, "85-3=", 2)
$d.Content.Find.Execute("14+0=", $true, $false, $false, $false, $false, $true, 1, $false, "13+45=", 2)
$d.Content.Find.Execute("12+50=", $true, $false, $false, $false, $false, $true, 1, $false, "29+7=", 2)
$d.Content.Find.Execute("49-46=", $true, $false, $false, $false, $false, $true, 1, $false, "44+31=", 2)
$d.Content.Find.Execute("99-32=", $true, $false, $false, $false, $false, $true, 1, $false, "85-32=", 2)
$d.Content.Find.Execute("20-16=", $true, $false, $false, $false, $false, $true, 1, $false, "92-32=", 2)
$d.Content.Find.Execute("45+26=", $true, $false, $false, $false, $false, $true, 1, $false, "43-39=", 2)
$d.Content.Find.Execute("2+21=", $true, $false, $false, $false, $false, $true, 1, $false, "0+32=", 2)
$d.Content.Find.Execute("97-20=", $true, $false, $false, $false, $false, $true, 1, $false, "46-4=", 2)
$d.Content.Find.Execute("4+9=", $true, $false, $false, $false, $false, $true, 1, $false, "67-67=", 2)
$d.Content.Find.Execute("73-30=", $true, $false, $false, $false, $false, $true, 1, $false, "5+61=", 2)
$d.Content.Find.Execute("21+40=", $true, $false, $false, $false, $false, $true, 1, $false, "97-4=", 2)
$d.Content.Find.Execute("14+44=", $true, $false, $false, $false, $false, $true, 1, $false, "3+30=", 2)
$d.Content.Find.Execute("5+22=", $true, $false, $false, $false, $false, $true, 1, $false, "45-32=", 2)
$d.Content.Find.Execute("27+57=", $true, $false, $false, $false, $false, $true, 1, $false, "68-33=", 2)
$d.Content.Find.Execute("42+5=", $true, $false, $false, $false, $false, $true, 1, $false, "31+11=", 2)
$d.Content.Find.Execute("18+69=", $true, $false, $false, $false, $false, $true, 1, $false, "21-2=", 2)
$d.Content.Find.Execute("54-38=", $true, $false, $false, $false, $false, $true, 1, $false, "11+45=", 2)
$d.Content.Find.Execute("48-26=", $true, $false, $false, $false, $false, $true, 1, $false, "61-43=", 2)
$d.Content.Find.Execute("70-62=", $true, $false, $false, $false, $false, $true, 1, $false, "88-31=", 2)
$d.Content.Find.Execute("83-41=", $true, $false, $false, $false, $false, $true, 1, $false, "0+9=", 2)
$d.Content.Find.Execute("79-42=", $true, $false, $false, $false, $false, $true, 1, $false, "79-24=", 2)
$d.Content.Find.Execute("44-41=", $true, $false, $false, $false, $false, $true, 1, $false, "6+71=", 2)
$d.Content.Find.Execute("69+19=", $true, $false, $false, $false, $false, $true, 1, $false, "98-92=", 2)
$d.Content.Find.Execute("29+37=", $true, $false, $false, $false, $false, $true, 1, $false, "22-21=", 2)
$d.Content.Find.Execute("87-72=", $true, $false, $false, $false, $false, $true, 1, $false, "10+33=", 2)
$d.Content.Find.Execute("12+72=", $true, $false, $false, $false, $false, $true, 1, $false, "86-69=", 2)
$d.Content.Find.Execute("32+3=", $true, $false, $false, $false, $false, $true, 1, $false, "48-5=", 2)
$d.Content.Find.Execute("30+35=", $true, $false, $false, $false, $false, $true, 1, $false, "75-16=", 2)
$d.Content.Find.Execute("11+60=", $true, $false, $false, $false, $false, $true, 1, $false, "43+24=", 2)
$d.Content.Find.Execute("5+42=", $true, $false, $false, $false, $false, $true, 1, $false, "80-14=", 2)
$d.Content.Find.Execute("2+15=", $true, $false, $false, $false, $false, $true, 1, $false, "7+52=", 2)
$d.Content.Find.Execute("60-50=", $true, $false, $false, $false, $false, $true, 1, $false, "77-13=", 2)
$d.Content.Find.Execute("35+36=", $true, $false, $false, $false, $false, $true, 1, $false, "53-6=", 2)
$d.Content.Find.Execute("92-76=", $true, $false, $false, $false, $false, $true, 1, $false, "91+7=", 2)
$d.Content.Find.Execute("21+21=", $true, $false, $false, $false, $false, $true, 1, $false, "74+21=", 2)
$d.Content.Find.Execute("20+13=", $true, $false, $false, $false, $false, $true, 1, $false, "31+65=", 2)
$d.Content.Find.Execute("58-10=", $true, $false, $false, $false, $false, $true, 1, $false, "46-42=", 2)
$d.Content.Find.Execute("58-1=", $true, $false, $false, $false, $false, $true, 1, $false, "20+49=", 2)
$d.Content.Find.Execute("56+1=", $true, $false, $false, $false, $false, $true, 1, $false, "44+47=", 2)
$d.Content.Find.Execute("89-52=", $true, $false, $false, $false, $false, $true, 1, $false, "16+13=", 2)
$d.Content.Find.Execute("20+34=", $true, $false, $false, $false, $false, $true, 1, $false, "62+8=", 2)
$d.Content.Find.Execute("87-62=", $true, $false, $false, $false, $false, $true, 1, $false, "19+18=", 2)
$d.Content.Find.Execute("35+18=", $true, $false, $false, $false, $false, $true, 1, $false, "32+27=", 2)
$d.Content.Find.Execute("41+36=", $true, $false, $false, $false, $false, $true, 1, $false, "53-44=", 2)
$d.Content.Find.Execute("63+31=", $true, $false, $false, $false, $false, $true, 1, $false, "14+27=", 2)
$d.Content.Find.Execute("20+67=", $true, $false, $false, $false, $false, $true, 1, $false, "71+3=", 2)
$d.Content.Find.Execute("91-51=", $true, $false, $false, $false, $false, $true, 1, $false, "95-30=", 2)
$d.Content.Find.Execute("83+15=", $true, $false, $false, $false, $false, $true, 1, $false, "22-13=", 2)
$d.Content.Find.Execute("35+4=", $true, $false, $false, $false, $false, $true, 1, $false, "39+14=", 2)
$d.Content.Find.Execute("36-4=", $true, $false, $false, $false, $false, $true, 1, $false, "88-66=", 2)
$d.Content.Find.Execute("42+36=", $true, $false, $false, $false, $false, $true, 1, $false, "39+49=", 2)
$d.Content.Find.Execute("74-70=", $true, $false, $false, $false, $false, $true, 1, $false, "98-92=", 2)
$d.Content.Find.Execute("93-9=", $true, $false, $false, $false, $false, $true, 1, $false, "41+3=", 2)
$d.Content.Find.Execute("46-5=", $true, $false, $false, $false, $false, $true, 1, $false, "10+48=", 2)
$d.Content.Find.Execute("61-30=", $true, $false, $false, $false, $false, $true, 1, $false, "43+32=", 2)
$d.Content.Find.Execute("7+1=", $true, $false, $false, $false, $false, $true, 1, $false, "54+11=", 2)
$d.Content.Find.Execute("66-16=", $true, $false, $false, $false, $false, $true, 1, $false, "31+23=", 2)
$d.Content.Find.Execute("37+29=", $true, $false, $false, $false, $false, $true, 1, $false, "0+78=", 2)
$d.Content.Find.Execute("61+5=", $true, $false, $false, $false, $false, $true, 1, $false, "10+2=", 2)
$d.Content.Find.Execute("12+78=", $true, $false, $false, $false, $false, $true, 1, $false, "74-52=", 2)
$d.Content.Find.Execute("81-6=", $true, $false, $false, $false, $false, $true, 1, $false, "89-88=", 2)
$d.Content.Find.Execute("51+21=", $true, $false, $false, $false, $false, $true, 1, $false, "18+22=", 2)
$d.Content.Find.Execute("17+73=", $true, $false, $false, $false, $false, $true, 1, $false, "96-18=", 2)
$d.Content.Find.Execute("83+12=", $true, $false, $false, $false, $false, $true, 1, $false, "14+26=", 2)
$d.Content.Find.Execute("33+34=", $true, $false, $false, $false, $false, $true, 1, $false, "61-7=", 2)
$d.Content.Find.Execute("44+55=", $true, $false, $false, $false, $false, $true, 1, $false, "40+19=", 2)
$d.Content.Find.Execute("39+36=", $true, $false, $false, $false, $false, $true, 1, $false, "97-78=", 2)
$d.Content.Find.Execute("60+39=", $true, $false, $false, $false, $false, $true, 1, $false, "93-54=", 2)
$d.Content.Find.Execute("40+50=", $true, $false, $false, $false, $false, $true, 1, $false, "1+94=", 2)
$d.Content.Find.Execute("2+87=", $true, $false, $false, $false, $false, $true, 1, $false, "64+2=", 2)
$d.Content.Find.Execute("70-41=", $true, $false, $false, $false, $false, $true, 1, $false, "58+10=", 2)
$d.Content.Find.Execute("14+74=", $true, $false, $false, $false, $false, $true, 1, $false, "76+16=", 2)
$d.Content.Find.Execute("12+4=", $true, $false, $false, $false, $false, $true, 1, $false, "65-24=", 2)
$d.Content.Find.Execute("50-27=", $true, $false, $false, $false, $false, $true, 1, $false, "24+61=", 2)
$d.Content.Find.Execute("86-5=", $true, $false, $false, $false, $false, $true, 1, $false, "18+55=", 2)
$d.Content.Find.Execute("42-39=", $true, $false, $false, $false, $false, $true, 1, $false, "26+26=", 2)
$d.Content.Find.Execute("39+4=", $true, $false, $false, $false, $false, $true, 1, $false, "73-61=", 2)
$d.Content.Find.Execute("72-24=", $true, $false, $false, $false, $false, $true, 1, $false, "37-13=", 2)
$d.Content.Find.Execute("23-1=", $true, $false, $false, $false, $false, $true, 1, $false, "53-19=", 2)
$d.Content.Find.Execute("52+39=", $true, $false, $false, $false, $false, $true, 1, $false, "62-30=", 2)
$d.Content.Find.Execute("38+8=", $true, $false, $false, $false, $false, $true, 1, $false, "82-67=", 2)
$d.Content.Find.Execute("59+18=", $true, $false, $false, $false, $false, $true, 1, $false, "66-36=", 2)
$d.Content.Find.Execute("71+21=", $true, $false, $false, $false, $false, $true, 1, $false, "48-27=", 2)
$d.Content.Find.Execute("11+3=", $true, $false, $false, $false, $false, $true, 1, $false, "70+0=", 2)
$d.Content.Find.Execute("79-36=", $true, $false, $false, $false, $false, $true, 1, $false, "23+74=", 2)
$d.Content.Find.Execute("84+14=", $true, $false, $false, $false, $false, $true, 1, $false, "76+9=", 2)
$d.Content.Find.Execute("36-16=", $true, $false, $false, $false, $false, $true, 1, $false, "23+38=", 2)
